$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-11 (columns A-G)
$data = @(
    @(26571, "Ana Freitas", "Vendas", "Problemas pessoais", 2, 45101, 9126.07),
    @(24627, "Agatha Carvalho", "Operações", "Outros", 8, 45090, 8311.030000000001),
    @(97429, "Isabelly Moraes", "Jurídico", "Problemas pessoais", 4, 45100, 4278.71),
    @(98068, "Sr. Otávio Costa", "Engenharia", "Doença", 1, 45085, 2887.44),
    @(33324, "Benício Castro", "Marketing", "Viagem de negócios", 5, 45093, 7499.03),
    @(90954, "Leandro Rezende", "Engenharia", "Doença", 2, 45088, 2988.57),
    @(74078, "Benício Cavalcanti", "P&D", "Outros", 6, 45096, 7410.12),
    @(34897, "Luigi Pereira", "P&D", "Doença", 5, 45089, 10906.86),
    @(35087, "Miguel Araújo", "Atendimento ao Cliente", "Doença", 2, 45086, 3905.96),
    @(13959, "Danilo Cardoso", "Operações", "Outros", 8, 45100, 11295.62)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
    $ws.Cells.Item($row, 7).Value = $values[6]
}

$wb.Save()
